$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comments (column H), entered in the same order the new unique remarks
#     were originally typed, so brand-new shared-string entries land in the
#     expected sequence ---
$ws.Range("H13").Value = "could not calibrate; tried four times"
$ws.Range("H17").Value = "calibration slightly more than 1! also I said her name. also said the noise bothered her. also played the gaze game with tami already "
$ws.Range("H6").Value = "he was a little distracted; may have said their name during the game to get their attention"
$ws.Range("H15").Value = "got very distracted by me; was not picking up her eyes because she was looking at me"
$ws.Range("H18").Value = "turned off the monitor at the beginning of the game, didn't finish"
$ws.Range("H19").Value = "headphones fell off for a sec, but other than that, it went well"
$ws.Range("H20").Value = "got a little frustrated towards the end, but went well"
$ws.Range("H21").Value = "all went well"
$ws.Range("H22").Value = "stopped playing at the end"

# --- Remaining H column updates (reuse of already-existing comments) ---
$ws.Range("H7").Value = "went well, but I didn't take his name tag off; also said the audio was too loud so I turned the volume down"
$ws.Range("H8").Value = "went well!"
$ws.Range("H9").Value = "went well!"
$ws.Range("H10").Value = "went well! When the noise condition started, she began to say the target words out loud"
$ws.Range("H11").Value = "went well!"
$ws.Range("H12").Value = "went well!"
$ws.Range("H14").Value = "went well!"
$ws.Range("H16").Value = "went well!"

# --- Fill in newly-collected raw data for subjects 12-22 (birthday, run date, gender) ---

# Row 12 - SPEED_ACC_NOISE_11
$ws.Range("B12").Value = "7/25/2017"
$ws.Range("C12").Value = "12/3/2012"
$ws.Range("D12").Value = "F"

# Row 13 - SPEED_ACC_NOISE_12
$ws.Range("B13").Value = "7/25/2017"
$ws.Range("C13").Value = "2/13/2014"
$ws.Range("D13").Value = "F"

# Row 14 - SPEED_ACC_NOISE_13
$ws.Range("B14").Value = "7/25/2017"
$ws.Range("C14").Value = "1/20/2014"
$ws.Range("D14").Value = "F"

# Row 15 - SPEED_ACC_NOISE_14
$ws.Range("B15").Value = "7/25/2017"
$ws.Range("C15").Value = "6/12/2012"
$ws.Range("D15").Value = "F"

# Row 16 - SPEED_ACC_NOISE_15
$ws.Range("B16").Value = "7/25/2017"
$ws.Range("C16").Value = "12/8/2013"
$ws.Range("D16").Value = "F"

# Row 17 - SPEED_ACC_NOISE_16
$ws.Range("B17").Value = "7/25/2017"
$ws.Range("C17").Value = "5/26/2012"
$ws.Range("D17").Value = "F"

# Row 18 - SPEED_ACC_NOISE_17
$ws.Range("B18").Value = "7/25/2017"
$ws.Range("C18").Value = "10/21/2013"
$ws.Range("D18").Value = "F"

# Row 19 - SPEED_ACC_NOISE_18
$ws.Range("B19").Value = "7/27/2017"
$ws.Range("C19").Value = "3/7/2013"
$ws.Range("D19").Value = "M"

# Row 20 - SPEED_ACC_NOISE_19
$ws.Range("B20").Value = "7/27/2017"
$ws.Range("C20").Value = "6/24/2014"
$ws.Range("D20").Value = "M"

# Row 21 - SPEED_ACC_NOISE_20
$ws.Range("B21").Value = "7/27/2017"
$ws.Range("C21").Value = "11/8/2011"
$ws.Range("D21").Value = "M"

# Row 22 - SPEED_ACC_NOISE_21
$ws.Range("B22").Value = "7/27/2017"
$ws.Range("C22").Value = "11/29/2012"
$ws.Range("D22").Value = "M"

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("A23").Select()
